# Add reference to figure 2a
#
# 1) Insert a superscript footnote-style marker "1" followed by a plain
#    space right after "...figure 2A from" (end of that paragraph).
# 2) Remove the oMathPara paragraph (the "@hirt2017" citation rendered
#    as an equation) that used to sit between "This is a replication of
#    figure 2A from" and "(fig.(...))".
# 3) Merge the (now adjacent) "This is a replication..." paragraph with
#    the "(fig.(...))" paragraph into a single paragraph.
# 4) Two unrelated text-only updates: the generation timestamp and the
#    git HEAD description near the end of the report.

$d = $word.ActiveDocument

# --- Step 1: insert "1" (superscript) + " " (plain) right after
#             "...figure 2A from", while it is still the end of its own
#             paragraph (nothing after it yet to merge formatting with).
$introRange = $d.Content
$null = $introRange.Find.Execute("This is a replication of figure 2A from", `
    $true, $false, $false, $false, $false, $true, 1, $false, "", 0)

$introRange.Collapse(0)  # wdCollapseEnd

$introRange.InsertAfter("1")
$oneRange = $d.Range($introRange.Start, $introRange.Start + 1)
$oneRange.Font.Superscript = $true

$spaceRange = $d.Range($oneRange.End, $oneRange.End)
$spaceRange.InsertAfter(" ")

# --- Step 2: delete the oMathPara paragraph -----------------------------
# It is the paragraph right after the one ending in
# "...figure 2A from1 " and its Range.Text renders as the (placeholder)
# math text "@???????2017".
$mathPara = $null
foreach ($p in $d.Paragraphs) {
    if ($p.Range.Text -like "*@*2017*") {
        $mathPara = $p
        break
    }
}
if ($mathPara -ne $null) {
    $mathPara.Range.Delete()
}

# --- Step 3: merge the now-adjacent "...figure 2A from1 " and
#             "(fig.(...))" paragraphs into one, by deleting the
#             paragraph mark between them.
$markRange = $d.Range($spaceRange.End, $spaceRange.End + 1)
$markRange.Delete()

# --- Step 4: plain text updates -----------------------------------------
$null = $d.Content.Find.Execute("This report was generated on 2021-11-30 10:21:14 using the following computational environment and dependencies:", `
    $true, $false, $false, $false, $false, $true, 1, $false, `
    "This report was generated on 2021-11-30 10:26:02 using the following computational environment and dependencies:", 2)

$null = $d.Content.Find.Execute("#> Head:     [810af1e] 2021-11-30: generate paper", `
    $true, $false, $false, $false, $false, $true, 1, $false, `
    "#> Head:     [0fd8cf5] 2021-11-30: Add figure reference", 2)
